# EarnedValue.xlsx edit script
# Implements: S-01030 bug fix row addition + S01024 -> S-01024 id fix +
# SUMIF range extension (131 -> 150) + two new "Horas insumidas" entries +
# hyperlink on the new story row + active-sheet/selection bookkeeping.

$wb = $excel.ActiveWorkbook
$wsEV = $wb.Worksheets.Item("Earned Value")
$wsHI = $wb.Worksheets.Item("Horas insumidas")

# ---------------------------------------------------------------------
# 1) "Horas insumidas": append the two new time entries for S-01030
#    (rows 136/137 sit in the gap between row 135 and the footer row 163,
#    so no row-insert is required there).
# ---------------------------------------------------------------------
$wsHI.Range("B136").Value2 = 40501
$wsHI.Range("C136").Value2 = "Sergio"
$wsHI.Range("D136").Value2 = "Se corrige el bug de prioridad media de validación de rangos para las métricas"
$wsHI.Range("E136").Value2 = "S-01030"
$wsHI.Range("F136").Value2 = 3

$wsHI.Range("B137").Value2 = 40502
$wsHI.Range("C137").Value2 = "Sergio"
$wsHI.Range("D137").Value2 = "Se corrige el bug de prioridad media de validación de rangos para las métricas"
$wsHI.Range("E137").Value2 = "S-01030"
$wsHI.Range("F137").Value2 = 6

# ---------------------------------------------------------------------
# 2) "Earned Value": fix the "S01024" id typo on row 22 (text of B22 is
#    unaffected, only A22 changes format to include the hyphen).
# ---------------------------------------------------------------------
$wsEV.Range("A22").Value2 = "S-01024"

# ---------------------------------------------------------------------
# 3) Insert a new row at position 23 (pushes the old rows 23-27 and the
#    trailing blank row 30 down by one) and populate it with the new
#    S-01030 story line.
# ---------------------------------------------------------------------
$wsEV.Rows.Item(23).Insert()

$wsEV.Range("A23").Value2 = "S-01030"
$wsEV.Range("B23").Value2 = "Arreglar todos los bugs de prioridad media/alta que figuran en el informe de avance"
$wsEV.Range("C23").Value2 = "Completada"
$wsEV.Range("D23").Value2 = 100
$wsEV.Range("E23").Value2 = 20
$wsEV.Range("F23").Value2 = 20
$wsEV.Range("G23").Formula = "=SUMIF('Horas insumidas'!`$E`$6:`$E`$150,A23,'Horas insumidas'!`$F`$6:`$F`$150)"
$wsEV.Range("H23").Formula = "=F23-G23"
$wsEV.Range("I23").Formula = "=F23-E23"
$wsEV.Range("J23").Formula = "=F23/E23"
$wsEV.Range("K23").Formula = "=F23/G23"

# Hyperlink on the new story's B23 cell.
$wsEV.Hyperlinks.Add($wsEV.Range("B23"), "https://www1.v1host.com/Team152/assetdetail.v1?oid=Story%3a1191", "", "", "https://www1.v1host.com/Team152/assetdetail.v1?oid=Story%3a1191")

# ---------------------------------------------------------------------
# 4) Update the SUMIF ranges on every other "Earned Value" data row
#    (2-22) so they look at 'Horas insumidas'!$E$6:$E$150 /
#    $F$6:$F$150 instead of the old $...$131 bound.
# ---------------------------------------------------------------------
for ($r = 2; $r -le 22; $r++) {
    $wsEV.Range("G$r").Formula = "=SUMIF('Horas insumidas'!`$E`$6:`$E`$150,A$r,'Horas insumidas'!`$F`$6:`$F`$150)"
}

# ---------------------------------------------------------------------
# 5) Selections / active sheet bookkeeping: end with "Horas insumidas"
#    selection at B138, then make "Earned Value" the active sheet with
#    K23 selected (matches the final workbook view state).
# ---------------------------------------------------------------------
$wsHI.Activate()
$wsHI.Range("B138").Select()

$wsEV.Activate()
$wsEV.Range("K23").Select()
